# Update 2p3. Added templates for formula student suspension, torque
# vectoring, four-wheel steering.
#
# For this workbook (sm_car_data_Tire_MFSwift.xlsx) the concrete change is:
# duplicate the existing "Tir_430_50R38" tire template sheet, place the
# copy right after it, rename the copy to "Tir_190_50R10", and update the
# two tire-specific cells (instance name + .tir file reference) on the new
# sheet so it describes the fsae190_50R10 tire instead of the truck tire.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Tir_430_50R38")

# Copy the template sheet, placing the new copy immediately after the
# source sheet (mirrors Excel's "Move or Copy... > Create a copy").
$source.Copy($null, $source)

# The freshly created copy is now the last sheet in the workbook.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Tir_190_50R10"

# Update the tire-specific cells: H5 (tirFile) first, then H3 (Instance),
# so they pick up shared-string slots in the same order as the source edit.
$newSheet.Range("H5").Value = "which('fsae190_50R10.tir')"
$newSheet.Range("H3").Value = "MFSwift_190_50R10"

# Make the newly added sheet the active sheet/tab, matching the saved
# workbook view after the edit.
$newSheet.Activate()
